$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name reorder (shared-string swaps reflected as displayed text) ---
$ws.Range("A66").Value = "Moldavia"
$ws.Range("A67").Value = "Azerbaiyan"
$ws.Range("A92").Value = "Grecia"
$ws.Range("A93").Value = "Croacia"
$ws.Range("A150").Value = "Reunion"
$ws.Range("A151").Value = "Trinidad yTobago"
$ws.Range("A152").Value = "Republica de Chipre"
$ws.Range("A153").Value = "Georgia"
$ws.Range("A174").Value = "Papua Nueva Guinea"
$ws.Range("A175").Value = "San Martin (Parte Holandesa)"
$ws.Range("A176").Value = "Burundi"

# --- Updated timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Agosto de 2020 a las 17:34"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Muertes hoy, Muertes) ---
$ws.Range("B4").Value = 6057928
$ws.Range("C4").Value = 11294
$ws.Range("D4").Value = 3350156
$ws.Range("E4").Value = 2522690
$ws.Range("G4").Value = 286
$ws.Range("H4").Value = 185082
$ws.Range("B6").Value = 3454513
$ws.Range("C6").Value = 69938
$ws.Range("D6").Value = 2640121
$ws.Range("E6").Value = 751723
$ws.Range("G6").Value = 975
$ws.Range("H6").Value = 62669
$ws.Range("B13").Value = 405972
$ws.Range("C13").Value = 1870
$ws.Range("D13").Value = 379452
$ws.Range("E13").Value = 15388
$ws.Range("G13").Value = 60
$ws.Range("H13").Value = 11132
$ws.Range("B66").Value = 35904
$ws.Range("C66").Value = 358
$ws.Range("D66").Value = 24156
$ws.Range("E66").Value = 10767
$ws.Range("G66").Value = 4
$ws.Range("H66").Value = 981
$ws.Range("B67").Value = 35844
$ws.Range("D67").Value = 33364
$ws.Range("E67").Value = 1956
$ws.Range("H67").Value = 524
$ws.Range("B92").Value = 9800
$ws.Range("C92").Value = 269
$ws.Range("D92").Value = 3804
$ws.Range("E92").Value = 5737
$ws.Range("G92").Value = 5
$ws.Range("H92").Value = 259
$ws.Range("B93").Value = 9549
$ws.Range("C93").Value = 357
$ws.Range("D93").Value = 6809
$ws.Range("E93").Value = 2560
$ws.Range("G93").Value = 3
$ws.Range("H93").Value = 180
$ws.Range("B117").Value = 3866
$ws.Range("C117").Value = 60
$ws.Range("D117").Value = 3222
$ws.Range("E117").Value = 552
$ws.Range("B150").Value = 1487
$ws.Range("C150").Value = 77
$ws.Range("D150").Value = 692
$ws.Range("E150").Value = 789
$ws.Range("H150").Value = 6
$ws.Range("B151").Value = 1476
$ws.Range("D151").Value = 594
$ws.Range("E151").Value = 867
$ws.Range("H151").Value = 15
$ws.Range("B152").Value = 1467
$ws.Range("C152").Value = 0
$ws.Range("D152").Value = 935
$ws.Range("E152").Value = 512
$ws.Range("H152").Value = 20
$ws.Range("B153").Value = 1455
$ws.Range("C153").Value = 8
$ws.Range("D153").Value = 1196
$ws.Range("E153").Value = 240
$ws.Range("H153").Value = 19
$ws.Range("B169").Value = 628
$ws.Range("C169").Value = 26
$ws.Range("E169").Value = 273
$ws.Range("B174").Value = 453
$ws.Range("C174").Value = 29
$ws.Range("D174").Value = 232
$ws.Range("E174").Value = 216
$ws.Range("G174").Value = 1
$ws.Range("H174").Value = 5
$ws.Range("B175").Value = 444
$ws.Range("C175").Value = 2
$ws.Range("D175").Value = 179
$ws.Range("E175").Value = 248
$ws.Range("H175").Value = 17
$ws.Range("B176").Value = 431
$ws.Range("D176").Value = 345
$ws.Range("E176").Value = 85
$ws.Range("H176").Value = 1
